# Generate Report for Handoff
# Adds two new handed-off files (3af2a31a-... and 5882ee1b-...) as new rows
# to the Overview / zh-cn / de-de sheets, mirroring the existing data pattern.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$neverHandback = "0001-01-01 00:00:00"
$handoffReason = "Include"
$mdExt = ".md"

$hyperUnderline = 2
$hyperColor = 15570276

# ---- File 1 : 3af2a31a-26ce-49ca-9f04-c85c29ccd5e1 ----
$uuid1 = "3af2a31a-26ce-49ca-9f04-c85c29ccd5e1"
$hash1 = "667968ef8cb52c9afaefb568608a97146d1d7c47"
$mdName1 = "$uuid1.md"
$xlfZh1 = "$uuid1.$hash1.zh-cn.xlf"
$xlfDe1 = "$uuid1.$hash1.de-de.xlf"

# ---- File 2 : 5882ee1b-fbb3-45f6-af16-96dc5f419f4c ----
$uuid2 = "5882ee1b-fbb3-45f6-af16-96dc5f419f4c"
$hash2 = "e5aa318102067d77b676c12946e186f5cdb58fa0"
$mdName2 = "$uuid2.md"
$xlfZh2 = "$uuid2.$hash2.zh-cn.xlf"
$xlfDe2 = "$uuid2.$hash2.de-de.xlf"

$dateOverview = "2016-30-11 09:30:08"
$dateZhCn = "2016-03-11 09:30:01"
$dateDeDe = "2016-03-11 09:30:08"

# Source GitHub blob URLs (same repo/path convention as the existing rows)
$srcUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/$hash1/e2e/$mdName1"
$srcUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/$hash2/e2e/$mdName2"

$handoffZhUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh1"
$handoffZhUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh2"

$handoffDeUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe1"
$handoffDeUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe2"

function Style-LikeHyperlink($range) {
    $range.Font.Underline = $hyperUnderline
    $range.Font.Color = $hyperColor
}

# =========================================================================
# Sheet "Overview" : columns A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
# =========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 6
$wsOverview.Cells.Item(6, 1).Value = $mdName1
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(6, 1), $srcUrl1, "", "", $mdName1)
Style-LikeHyperlink($wsOverview.Cells.Item(6, 1))
$wsOverview.Cells.Item(6, 2).Value = $statusReady
$wsOverview.Cells.Item(6, 3).Value = $statusReady
$wsOverview.Cells.Item(6, 4).Value = $dateOverview

# Row 7
$wsOverview.Cells.Item(7, 1).Value = $mdName2
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(7, 1), $srcUrl2, "", "", $mdName2)
Style-LikeHyperlink($wsOverview.Cells.Item(7, 1))
$wsOverview.Cells.Item(7, 2).Value = $statusReady
$wsOverview.Cells.Item(7, 3).Value = $statusReady
$wsOverview.Cells.Item(7, 4).Value = $dateOverview

# =========================================================================
# Sheet "zh-cn" : A=Source File Name, B=File Extension, C=Status,
#                 D=Latest Handoff File, E=Latest Handoff Datetime,
#                 H=Latest Handback DateTime, I=Handoff Reason
# =========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 6
$wsZhCn.Cells.Item(6, 1).Value = $mdName1
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(6, 1), $srcUrl1, "", "", $mdName1)
Style-LikeHyperlink($wsZhCn.Cells.Item(6, 1))
$wsZhCn.Cells.Item(6, 2).Value = $mdExt
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(6, 2), $srcUrl1, "", "", $mdExt)
Style-LikeHyperlink($wsZhCn.Cells.Item(6, 2))
$wsZhCn.Cells.Item(6, 3).Value = $statusReady
$wsZhCn.Cells.Item(6, 4).Value = $xlfZh1
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(6, 4), $handoffZhUrl1, "", "", $xlfZh1)
Style-LikeHyperlink($wsZhCn.Cells.Item(6, 4))
$wsZhCn.Cells.Item(6, 5).Value = $dateZhCn
$wsZhCn.Cells.Item(6, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(6, 8).Value = $neverHandback
$wsZhCn.Cells.Item(6, 9).Value = $handoffReason

# Row 7
$wsZhCn.Cells.Item(7, 1).Value = $mdName2
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(7, 1), $srcUrl2, "", "", $mdName2)
Style-LikeHyperlink($wsZhCn.Cells.Item(7, 1))
$wsZhCn.Cells.Item(7, 2).Value = $mdExt
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(7, 2), $srcUrl2, "", "", $mdExt)
Style-LikeHyperlink($wsZhCn.Cells.Item(7, 2))
$wsZhCn.Cells.Item(7, 3).Value = $statusReady
$wsZhCn.Cells.Item(7, 4).Value = $xlfZh2
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(7, 4), $handoffZhUrl2, "", "", $xlfZh2)
Style-LikeHyperlink($wsZhCn.Cells.Item(7, 4))
$wsZhCn.Cells.Item(7, 5).Value = $dateZhCn
$wsZhCn.Cells.Item(7, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(7, 8).Value = $neverHandback
$wsZhCn.Cells.Item(7, 9).Value = $handoffReason

# =========================================================================
# Sheet "de-de" : same column layout as zh-cn
# =========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 6
$wsDeDe.Cells.Item(6, 1).Value = $mdName1
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(6, 1), $srcUrl1, "", "", $mdName1)
Style-LikeHyperlink($wsDeDe.Cells.Item(6, 1))
$wsDeDe.Cells.Item(6, 2).Value = $mdExt
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(6, 2), $srcUrl1, "", "", $mdExt)
Style-LikeHyperlink($wsDeDe.Cells.Item(6, 2))
$wsDeDe.Cells.Item(6, 3).Value = $statusReady
$wsDeDe.Cells.Item(6, 4).Value = $xlfDe1
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(6, 4), $handoffDeUrl1, "", "", $xlfDe1)
Style-LikeHyperlink($wsDeDe.Cells.Item(6, 4))
$wsDeDe.Cells.Item(6, 5).Value = $dateDeDe
$wsDeDe.Cells.Item(6, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(6, 8).Value = $neverHandback
$wsDeDe.Cells.Item(6, 9).Value = $handoffReason

# Row 7
$wsDeDe.Cells.Item(7, 1).Value = $mdName2
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(7, 1), $srcUrl2, "", "", $mdName2)
Style-LikeHyperlink($wsDeDe.Cells.Item(7, 1))
$wsDeDe.Cells.Item(7, 2).Value = $mdExt
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(7, 2), $srcUrl2, "", "", $mdExt)
Style-LikeHyperlink($wsDeDe.Cells.Item(7, 2))
$wsDeDe.Cells.Item(7, 3).Value = $statusReady
$wsDeDe.Cells.Item(7, 4).Value = $xlfDe2
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(7, 4), $handoffDeUrl2, "", "", $xlfDe2)
Style-LikeHyperlink($wsDeDe.Cells.Item(7, 4))
$wsDeDe.Cells.Item(7, 5).Value = $dateDeDe
$wsDeDe.Cells.Item(7, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(7, 8).Value = $neverHandback
$wsDeDe.Cells.Item(7, 9).Value = $handoffReason
